$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7
$ws.Range("B3").Value = 117940
$ws.Range("B4").Value = 60
$ws.Range("B5").Value = 60
